# Update res_bus/vm_pu.xlsx results for the "case with 380 kV" run: rows 2-25
# (bus indices 0-23), columns C,D,E,F,J,K,L,M,N get new per-unit voltage values
# from the re-run power flow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "C2" = 1.077097194815221
    "D2" = 1.08692038159443
    "E2" = 1.081559366048155
    "F2" = 1.094771762609404
    "J2" = 1.081993357721314
    "K2" = 1.089576980966634
    "L2" = 1.084229872939023
    "M2" = 1.097408259026152
    "N2" = 1.030231680001454
    "C3" = 1.079629872643857
    "D3" = 1.089394839400521
    "E3" = 1.083873224589527
    "F3" = 1.097286405895668
    "J3" = 1.084178798154688
    "K3" = 1.091866945882647
    "L3" = 1.086358584309257
    "M3" = 1.099739819178847
    "N3" = 1.031019165620976
    "C4" = 1.081261425383544
    "D4" = 1.090989084607025
    "E4" = 1.085363605105732
    "F4" = 1.098906681326331
    "J4" = 1.085585655038974
    "K4" = 1.093341504606985
    "L4" = 1.087728826230067
    "M4" = 1.101241321932018
    "N4" = 1.031524763764052
    "C5" = 1.081945641130516
    "D5" = 1.091657701279923
    "E5" = 1.085988567002672
    "F5" = 1.099586248201469
    "J5" = 1.086175400095348
    "K5" = 1.093959729160351
    "L5" = 1.08830319924295
    "M5" = 1.101870881242414
    "N5" = 1.031736384892076
    "C6" = 1.082060426328191
    "D6" = 1.091769872243303
    "E6" = 1.086093408635332
    "F6" = 1.099700258182677
    "J6" = 1.086274322569442
    "K6" = 1.094063434544671
    "L6" = 1.088399541871633
    "M6" = 1.101976490236201
    "N6" = 1.031771862775023
    "C7" = 1.081270574479412
    "D7" = 1.090998024936303
    "E7" = 1.085371962082601
    "F7" = 1.09891576794419
    "J7" = 1.085593541854631
    "K7" = 1.093349771895773
    "L7" = 1.087736507563335
    "M7" = 1.101249740646277
    "N7" = 1.031527595099062
    "C8" = 1.077954657831719
    "D8" = 1.087758091969407
    "E8" = 1.082342788236642
    "F8" = 1.095623048461126
    "J8" = 1.082733468584535
    "K8" = 1.090352403287819
    "L8" = 1.084950792022394
    "M8" = 1.098197732932352
    "N8" = 1.030498643121221
    "C9" = 1.072053773455259
    "D9" = 1.081993998353002
    "E9" = 1.076950624113911
    "F9" = 1.089766126209511
    "J9" = 1.077636040674979
    "K9" = 1.085013493287195
    "L9" = 1.079985158096745
    "M9" = 1.092762741635
    "N9" = 1.028654497441378
    "C10" = 1.068077840655972
    "D10" = 1.078111355415872
    "E10" = 1.073316488630716
    "F10" = 1.085821691539815
    "J10" = 1.074196275556258
    "K10" = 1.081412947653394
    "L10" = 1.076633851925781
    "M10" = 1.089098243217101
    "N10" = 1.027403246845007
    "C11" = 1.066345544997929
    "D11" = 1.076419985771268
    "E11" = 1.071732899355284
    "F11" = 1.08410358040452
    "J11" = 1.072696360970624
    "K11" = 1.079843441168943
    "L11" = 1.075172400775277
    "M11" = 1.087501061105504
    "N11" = 1.026856044875795
    "C12" = 1.065700425458805
    "D12" = 1.075790150631391
    "E12" = 1.071143128613404
    "F12" = 1.083463813866094
    "J12" = 1.072137598405156
    "K12" = 1.079258831294015
    "L12" = 1.074627950888349
    "M12" = 1.086906173083126
    "N12" = 1.026651958764889
    "C13" = 1.065838882266835
    "D13" = 1.075925325096328
    "E13" = 1.071269707712811
    "F13" = 1.083601118588459
    "J13" = 1.072257529400075
    "K13" = 1.079384306554309
    "L13" = 1.07474481057869
    "M13" = 1.087033852969397
    "N13" = 1.026695773845097
    "C14" = 1.06629225359982
    "D14" = 1.07636795608323
    "E14" = 1.071684180769943
    "F14" = 1.084050729755418
    "J14" = 1.072650206999685
    "K14" = 1.079795150616917
    "L14" = 1.07512742936232
    "M14" = 1.087451920925396
    "N14" = 1.026839192118332
    "C15" = 1.066571367774509
    "D15" = 1.076640463816594
    "E15" = 1.071939343556082
    "F15" = 1.084327538398556
    "J15" = 1.072891931335021
    "K15" = 1.080048068147762
    "L15" = 1.075362959575611
    "M15" = 1.087709289556107
    "N15" = 1.026927446157972
    "C16" = 1.068192576942075
    "D16" = 1.078223387008181
    "E16" = 1.073421371194165
    "F16" = 1.085935498133602
    "J16" = 1.074295594700405
    "K16" = 1.081516885796701
    "L16" = 1.076730621886968
    "M16" = 1.089204018390467
    "N16" = 1.02743944738014
    "C17" = 1.069206613332418
    "D17" = 1.079213551775578
    "E17" = 1.074348296271925
    "F17" = 1.086941371166384
    "J17" = 1.075173233824029
    "K17" = 1.082435400445557
    "L17" = 1.077585722482316
    "M17" = 1.090138790368429
    "N17" = 1.027759152391226
    "C18" = 1.069797056351689
    "D18" = 1.079790121905401
    "E18" = 1.074887996393727
    "F18" = 1.087527105126905
    "J18" = 1.075684139027547
    "K18" = 1.082970150395693
    "L18" = 1.07808349668462
    "M18" = 1.090683024565268
    "N18" = 1.027945111006379
    "C19" = 1.069998209770691
    "D19" = 1.079986553595487
    "E19" = 1.075071859074925
    "F19" = 1.087726661713932
    "J19" = 1.075858175336502
    "K19" = 1.083152317668115
    "L19" = 1.078253058215992
    "M19" = 1.090868426015081
    "N19" = 1.028008430481286
    "C20" = 1.069097923520201
    "D20" = 1.079107417864151
    "E20" = 1.074248945633897
    "F20" = 1.086833551759722
    "J20" = 1.075079175982125
    "K20" = 1.082336956776089
    "L20" = 1.07749408121301
    "M20" = 1.090038602278287
    "N20" = 1.02772490498434
    "C21" = 1.066158793583258
    "D21" = 1.076237656525636
    "E21" = 1.071562172163621
    "F21" = 1.083918374670248
    "J21" = 1.072534618599805
    "K21" = 1.07967421260026
    "L21" = 1.075014802312156
    "M21" = 1.087328855607379
    "N21" = 1.026796982106493
    "C22" = 1.064301165077551
    "D22" = 1.074424122510584
    "E22" = 1.069863865158009
    "F22" = 1.082076294799279
    "J22" = 1.070925310595792
    "K22" = 1.077990607768553
    "L22" = 1.073446686672099
    "M22" = 1.085615708601241
    "N22" = 1.026208742761088
    "C23" = 1.065286867658084
    "D23" = 1.075386403314035
    "E23" = 1.070765044081621
    "F23" = 1.083053707678578
    "J23" = 1.071779348693187
    "K23" = 1.078884031521249
    "L23" = 1.074278873237088
    "M23" = 1.086524792474834
    "N23" = 1.026521042759848
    "C24" = 1.06914703893907
    "D24" = 1.079155378223979
    "E24" = 1.074293840873249
    "F24" = 1.086882273717742
    "J24" = 1.075121679780781
    "K24" = 1.082381442332871
    "L24" = 1.07753549301921
    "M24" = 1.090083876058647
    "N24" = 1.027740381518992
    "C25" = 1.073586464254427
    "D25" = 1.083490967031528
    "E25" = 1.078351356378977
    "F25" = 1.091287074335803
    "J25" = 1.078960954643529
    "K25" = 1.08640079251871
    "L25" = 1.0812759036997992
    "M25" = 1.094174859860206
    "N25" = 1.029135026410763
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}

Write-Host "Applied $($data.Count) cell updates"
